# Regenerate the "K" column (column G, formerly "Strike#") with recalculated
# strike-count values for each trade row (rows 2-76 of Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @(
    1, 3, 2, 0, 2, 0, 0, 0, 1, 3, 1, 0, 2, 0, 0, 2, 0, 1, 0, 1,
    1, 2, 1, 1, 1, 0, 2, 3, 1, 0, 2, 0, 0, 2, 2, 2, 0, 3, 1, 0,
    1, 0, 0, 1, 1, 2, 2, 3, 3, 3, 2, 2, 1, 0, 3, 1, 2, 0, 0, 0,
    2, 2, 1, 1, 1, 2, 2, 3, 1, 2, 2, 0, 1, 1, 1
)

for ($i = 0; $i -lt $kValues.Length; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 7).Value = $kValues[$i]
}
